# Update the "Förändrad" (Changed) date column C for all data rows (2-173)
# from 45175 (2023-09-06) to 45183 (2023-09-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C173").Value = 45183
